# Restructure the "LKT Events" sheet (sheet1) so that the request data
# (Long name) is isolated into its own column, separate from the
# logic/description columns.
#
# This inserts a new column C ("Long name") between the existing
# "Short label" (B) and "Description in text" (old C, now D) columns,
# and populates it with the short/long-name mapping for each event row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LKT Events")

# Insert a new, blank column before column C. This shifts the old
# "Description in text" and "HED tags" columns one position to the right.
$ws.Columns("C").Insert()

# Match the look of the surrounding data columns (center/wrap like the
# rest of the table) and give it roughly the same width as column B.
$newCol = $ws.Columns("C")
$newCol.ColumnWidth = $ws.Columns("B").ColumnWidth()
$newCol.VerticalAlignment = -4108
$newCol.WrapText = $True

# Header
$ws.Range("C1").Value = "Long name"

# Row data - the "long name" associated with each event's short label.
$ws.Range("C2").Value = "PerturbCarToLeft"
$ws.Range("C3").Value = "PerturbCarToRight"
$ws.Range("C5").Value = "DriverStopsCorrecting"
$ws.Range("C4").Value = "DriverStartsToCorrect"

# Leave the selection on the newly added cell, as in the authored edit.
$ws.Range("C4").Select() | Out-Null
